# Updated final report outline
#
# 1) The three name paragraphs (Huiting / Xingliang / Seoho) were wrapped
#    in <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/>
#    markers. Those transient proofing markers are not part of a Range's
#    WordOpenXML projection, so reinserting a paragraph's own
#    WordOpenXML into itself is a clean, attribute-preserving way to drop
#    them without disturbing anything else about the paragraph.
$d = $word.ActiveDocument

# The Huiting / Xingliang / Seoho runs are the first three paragraphs in
# the document.
for ($i = 1; $i -le 3; $i++) {
    $p = $d.Paragraphs($i)
    $r = $p.Range
    $r.InsertXML($r.WordOpenXML) | Out-Null
}

# 2) Append ": what we did (converting into the required format)" right
#    after the "Conversion and EPOS" run, in the same run formatting, and
#    move the _GoBack bookmark so it now sits at the end of that
#    paragraph instead of at the end of the "First set" paragraph.
$rng = $d.Content
$rng.Find.Execute("Conversion and EPOS") | Out-Null
$convPara = $rng.Paragraphs(1)
$convRange = $convPara.Range
$convXml = $convRange.WordOpenXML
$anchor = '<w:t>Conversion and EPOS</w:t></w:r>'
$addition = '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="4472C4" w:themeColor="accent5"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:eastAsia="ko-KR"/></w:rPr><w:t>: what we did (converting into the required format)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$convXml = $convXml -replace [regex]::Escape($anchor), ($anchor + $addition)
$convRange.InsertXML($convXml) | Out-Null

# Remove the _GoBack bookmark pair from the end of the "First set"
# paragraph now that it has moved.
$rng2 = $d.Content
$rng2.Find.Execute("First set") | Out-Null
$fsPara = $rng2.Paragraphs(1)
$fsRange = $fsPara.Range
$fsXml = $fsRange.WordOpenXML
$fsXml = $fsXml -replace '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>', ''
$fsRange.InsertXML($fsXml) | Out-Null
